$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 5 ("文档管理" / Document Management) - name/uri updated to the new
# integration-based document endpoint.
# ---------------------------------------------------------------------
$ws.Range("F5").Value2 = "zero.document"
$ws.Range("I5").Value2 = "/integration/document"

# ---------------------------------------------------------------------
# Row 6 ("集成管理" / Integration Management) - brand-new top level
# DEV-MENU container row.
# ---------------------------------------------------------------------
$ws.Range("A6").Value2 = "24719b19-ea7f-46a2-a808-65c1839ca001"
$ws.Range("C6").Value2 = "DEV-MENU"
$ws.Range("D6").Value2 = 70000
$ws.Range("E6").Value2 = 1
$ws.Range("F6").Value2 = "develop.integration"
$ws.Range("G6").Value2 = "集成管理"
$ws.Range("H6").Value2 = "api"
$ws.Range("I6").Value2 = "EXPAND"

# ---------------------------------------------------------------------
# Row 7 ("FTP配置") - child of row 6, parentId references A$6.
# ---------------------------------------------------------------------
$ws.Range("A7").Value2 = "af4567a8-9551-44a7-b26a-4b723dbd9d6f"
$ws.Range("B7").Formula = "=A`$6"
$ws.Range("C7").Value2 = "DEV-MENU"
$ws.Range("D7").Value2 = 1005
$ws.Range("E7").Value2 = 2
$ws.Range("F7").Value2 = "develop.integration.ftp"
$ws.Range("G7").Value2 = "FTP配置"
$ws.Range("H7").Value2 = "folder-open"
$ws.Range("I7").Value2 = "/integration/ftp"

# ---------------------------------------------------------------------
# Row 8 ("邮件配置") - child of row 6.
# ---------------------------------------------------------------------
$ws.Range("A8").Value2 = "c7627312-6912-47a6-b1df-f2a334ddf450"
$ws.Range("B8").Formula = "=A`$6"
$ws.Range("C8").Value2 = "DEV-MENU"
$ws.Range("D8").Value2 = 1005
$ws.Range("E8").Value2 = 2
$ws.Range("F8").Value2 = "develop.integration.email"
$ws.Range("G8").Value2 = "邮件配置"
$ws.Range("H8").Value2 = "mail"
$ws.Range("I8").Value2 = "/integration/email"

# ---------------------------------------------------------------------
# Row 9 ("短信配置") - child of row 6.
# ---------------------------------------------------------------------
$ws.Range("A9").Value2 = "4801ae80-8f45-48ae-aec1-a897723958f4"
$ws.Range("B9").Formula = "=A`$6"
$ws.Range("C9").Value2 = "DEV-MENU"
$ws.Range("D9").Value2 = 1005
$ws.Range("E9").Value2 = 2
$ws.Range("F9").Value2 = "develop.integration.sms"
$ws.Range("G9").Value2 = "短信配置"
$ws.Range("H9").Value2 = "message"
$ws.Range("I9").Value2 = "/integration/sms"

# ---------------------------------------------------------------------
# Formatting: copy existing row-5 look & feel down onto the new rows so
# borders/fonts/fills/number formats line up with the rest of the table.
# ---------------------------------------------------------------------
$ws.Range("A5:K5").Copy()
$ws.Range("A6:K6").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:K5").Copy()
$ws.Range("A7:K7").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:K5").Copy()
$ws.Range("A8:K8").PasteSpecial(-4122) | Out-Null
$ws.Range("A5:K5").Copy()
$ws.Range("A9:K9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# D6 keeps the plain "order" look used by D5 (copied above); D7:D9 get the
# centred variant used for the nested FTP/email/SMS rows.
$ws.Range("D5").Copy()
$ws.Range("D7:D9").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# I6 ("EXPAND") is flagged in red to highlight the new development center
# entry point, re-using the existing bold red header font + thin border.
$ws.Range("F5").Copy()
$ws.Range("I6").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("I6").Font.Color = 255
$ws.Range("I6").Font.Size = 16
$ws.Range("I6").Font.Name = "等线"
$ws.Range("I6").HorizontalAlignment = -4131

Write-Host "Applied Development Center (Integration) menu rows."
